$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency Price (D) and Volume(1h) (E) columns.
# Numeric-looking Price strings must be force-formatted as Text first so
# Excel keeps them as literal strings (preserving trailing zeros / exact digits)
# instead of silently coercing them to numbers.
$updates = @(
    @{ Cell = 'D2'; Value = '68.854.09' },
    @{ Cell = 'E2'; Value = '  +0.42%  ' },
    @{ Cell = 'D3'; Value = '2.472.93' },
    @{ Cell = 'E3'; Value = '  +0.49%  ' },
    @{ Cell = 'E4'; Value = '  -0.02%  ' },
    @{ Cell = 'D5'; Value = '561.23' },
    @{ Cell = 'E5'; Value = '  +0.19%  ' },
    @{ Cell = 'D6'; Value = '164.07' },
    @{ Cell = 'E6'; Value = '  -0.31%  ' },
    @{ Cell = 'D8'; Value = '0.514' },
    @{ Cell = 'E8'; Value = '  +1.88%  ' },
    @{ Cell = 'D9'; Value = '0.157' },
    @{ Cell = 'E9'; Value = '  +3.81%  ' },
    @{ Cell = 'E11'; Value = '  -1.64%  ' },
    @{ Cell = 'D12'; Value = '4.85' },
    @{ Cell = 'E12'; Value = '  +0.41%  ' },
    @{ Cell = 'D13'; Value = '68.724.92' },
    @{ Cell = 'E13'; Value = '  +0.43%  ' },
    @{ Cell = 'D14'; Value = '0.0000171' },
    @{ Cell = 'E14'; Value = '  +0.29%  ' },
    @{ Cell = 'D15'; Value = '23.61' },
    @{ Cell = 'E15'; Value = '  +1.07%  ' },
    @{ Cell = 'D16'; Value = '10.65' },
    @{ Cell = 'E16'; Value = '  -3.35%  ' },
    @{ Cell = 'D17'; Value = '338.12' },
    @{ Cell = 'E17'; Value = '  -2.08%  ' },
    @{ Cell = 'D18'; Value = '6.94' },
    @{ Cell = 'E18'; Value = '  -3.60%  ' },
    @{ Cell = 'D19'; Value = '3.80' },
    @{ Cell = 'E19'; Value = '  +0.26%  ' },
    @{ Cell = 'D20'; Value = '1.89' },
    @{ Cell = 'E20'; Value = '  +0.67%  ' },
    @{ Cell = 'E21'; Value = '  +0.01%  ' },
    @{ Cell = 'D22'; Value = '66.73' },
    @{ Cell = 'E22'; Value = '  -1.80%  ' },
    @{ Cell = 'D23'; Value = '3.67' },
    @{ Cell = 'E23'; Value = '  -1.58%  ' },
    @{ Cell = 'D24'; Value = '8.26' },
    @{ Cell = 'E24'; Value = '  +1.24%  ' },
    @{ Cell = 'D25'; Value = '0.0₃0826' },
    @{ Cell = 'E25'; Value = '  -1.25%  ' },
    @{ Cell = 'D26'; Value = '7.23' },
    @{ Cell = 'E26'; Value = '  -0.18%  ' },
    @{ Cell = 'D27'; Value = '0.999' },
    @{ Cell = 'E27'; Value = '  +0.02%  ' },
    @{ Cell = 'D28'; Value = '430.61' },
    @{ Cell = 'E28'; Value = '  -1.07%  ' },
    @{ Cell = 'D29'; Value = '1.15' },
    @{ Cell = 'E29'; Value = '  -1.81%  ' },
    @{ Cell = 'D30'; Value = '1.63' },
    @{ Cell = 'E30'; Value = '  -2.79%  ' },
    @{ Cell = 'D31'; Value = '160.20' },
    @{ Cell = 'E31'; Value = '  +1.68%  ' },
    @{ Cell = 'E32'; Value = '  +0.00%  ' },
    @{ Cell = 'E33'; Value = '  -0.06%  ' },
    @{ Cell = 'E34'; Value = '  -1.89%  ' },
    @{ Cell = 'D35'; Value = '17.88' },
    @{ Cell = 'E35'; Value = '  -0.28%  ' },
    @{ Cell = 'D36'; Value = '4.46' },
    @{ Cell = 'E36'; Value = '  -0.42%  ' },
    @{ Cell = 'D37'; Value = '0.298' },
    @{ Cell = 'E37'; Value = '  -2.76%  ' },
    @{ Cell = 'D38'; Value = '1.48' },
    @{ Cell = 'E38'; Value = '  -3.03%  ' },
    @{ Cell = 'D39'; Value = '1.08' },
    @{ Cell = 'E39'; Value = '  -2.14%  ' },
    @{ Cell = 'D40'; Value = '2.07' },
    @{ Cell = 'E40'; Value = '  -1.41%  ' },
    @{ Cell = 'D41'; Value = '3.39' },
    @{ Cell = 'E41'; Value = '  +0.82%  ' },
    @{ Cell = 'D42'; Value = '130.62' },
    @{ Cell = 'E42'; Value = '  -3.28%  ' },
    @{ Cell = 'E43'; Value = '  +0.32%  ' },
    @{ Cell = 'D44'; Value = '0.486' },
    @{ Cell = 'E44'; Value = '  +0.31%  ' },
    @{ Cell = 'D45'; Value = '0.566' },
    @{ Cell = 'E45'; Value = '  +0.33%  ' },
    @{ Cell = 'D46'; Value = '0.0918' },
    @{ Cell = 'E46'; Value = '  +0.98%  ' },
    @{ Cell = 'E47'; Value = '  +0.00%  ' },
    @{ Cell = 'D48'; Value = '1.39' },
    @{ Cell = 'E48'; Value = '  -2.61%  ' },
    @{ Cell = 'D49'; Value = '5.01' },
    @{ Cell = 'E49'; Value = '  -7.02%  ' },
    @{ Cell = 'D50'; Value = '16.92' },
    @{ Cell = 'E50'; Value = '  -4.08%  ' },
    @{ Cell = 'D51'; Value = '0.0₆0206' },
    @{ Cell = 'E51'; Value = '  -8.95%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $looksNumeric = $u.Value -match '^[+-]?\d+(\.\d+)?$'
    if ($looksNumeric) {
        # Force text storage so "3.80" / "160.20" keep their exact digits
        # rather than becoming the numbers 3.8 / 160.2.
        $cell.NumberFormat = "@"
    }
    $cell.Value = $u.Value
}
